# Ready for manufacturing (V2.0)
# Adds "Supplier 1" / "Supplier Part Number 1" columns (I, J) to the BOM sheet,
# and refreshes the "LED W" (D513/D514/D515) row with updated part info.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-CellLikeNeighbor {
    param($targetCell, $formatSourceCell, $value)
    $targetCell.Value = $value
    $formatSourceCell.Copy() | Out-Null
    $targetCell.PasteSpecial($xlPasteFormats) | Out-Null
}

# --- New column widths for I:J (approx. 14.43 chars, matching F:H band) ---
$ws.Columns.Item(9).ColumnWidth = 13.59
$ws.Columns.Item(10).ColumnWidth = 13.59

# --- Header row (row 1): new headers, styled like the rest of the header row (H1) ---
Set-CellLikeNeighbor $ws.Cells.Item(1,9) $ws.Cells.Item(1,8) "Supplier 1"
Set-CellLikeNeighbor $ws.Cells.Item(1,10) $ws.Cells.Item(1,8) "Supplier Part Number 1"

# --- Update the LED W (D513, D514, D515) row 5 part info, preserving existing cell styles ---
# (use the same column one row down as the format donor, since row 5's own cells are
# the ones being overwritten)
Set-CellLikeNeighbor $ws.Cells.Item(5,5) $ws.Cells.Item(6,5) "LED XLAMP WARM WHT 2700K SMD"
Set-CellLikeNeighbor $ws.Cells.Item(5,7) $ws.Cells.Item(6,7) "XTEAWT-E0-0000-00000HEE8"
Set-CellLikeNeighbor $ws.Cells.Item(5,8) $ws.Cells.Item(6,8) "D_LED_W_HP_3"

# --- Per-row Supplier 1 / Supplier Part Number 1 values (styled like column H of the same row) ---
$supplierData = @{
    2  = @("Digi-Key", "490-10450-1-ND")
    3  = @("LCSC", "C293627")
    4  = @("Digi-Key", "404-1273-1-ND")
    5  = @("Digi-Key", "XTEAWT-E0-0000-00000HEE8TR-ND")
    6  = @("Digi-Key", "475-3442-1-ND")
    7  = @("LCSC", "C710226")
    8  = @("LCSC", "C524044")
    9  = @("Mouser", "998-MIC2841AYMTTR")
    10 = @("LCSC", "C842736")
    11 = @("LCSC", "C50299")
    12 = @("Amazon", "B077VQTB6Q")
    13 = @("LCSC", "C294565")
    14 = @("Digi-Key", "153-1135-ND")
    15 = @("Digi-Key", "RMCF0402FT2K00CT-ND")
    16 = @("LCSC", "C104939")
    17 = @("Digi-Key", "311-10.0KLRCT-ND")
    18 = @("Digi-Key", "RMCF0402ZT0R00CT-ND")
    19 = @("Digi-Key", "RMCF0402FT16K0CT-ND")
}

foreach ($row in 2..19) {
    $vals = $supplierData[$row]
    $refCell = $ws.Cells.Item($row, 8)

    Set-CellLikeNeighbor $ws.Cells.Item($row, 9) $refCell $vals[0]
    Set-CellLikeNeighbor $ws.Cells.Item($row, 10) $refCell $vals[1]
}
